# Apply the metadata refresh for the Alvearie salary-indicator StructureDefinition export.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# Version bump
$meta.Range("B3").Value = "6.0.0"

# Regeneration timestamp
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was previously blank
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 used to be a duplicate "Contact" row with no real display text;
# turn it into the Jurisdiction row and drop the extra duplicate row (row 11).
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"
$meta.Rows.Item(11).Delete()

# The root Extension slice row on the Elements sheet now carries the
# specific extension's short/definition text instead of the generic ones.
$elements.Range("K2").Value = "Salary Indicator"
$elements.Range("L2").Value = "Indicator of whether the employee status is salaried"
